$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.239.38"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.575.14"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'208.03"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'22.21"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.799.07"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "1.583.02"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "27.251.01"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "'62.36"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "'214.59"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "'7.34"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "0.0₃0686"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'4.12"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'152.02"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "'6.68"
$ws.Range("E26").Value = "  -4.53%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").Value = "'0.0464"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").Value = "1.406.59"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").Value = "'0.0165"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("D43").Value = "'1.81"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").Value = "'5.36"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "'63.77"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "1.711.41"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'86.00"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "0.0₇0982"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E51").Value = "  -0.12%  "
